# SC-IMT ET TT Overview - apply commit:
#   "update the region & gender"
#   "add TT-1st Mile two Vice-Chair"
#
# Changes:
#  1. Header G1 "Vice-Chair" -> "Co-Chair/Vice-Chair*"
#  2. G12 "Eugene Burger" -> "Eugene Burger*" (footnote marker added)
#  3. TT-First-Mile (row 11-13 block) gets a 2nd Vice-Chair: insert a new row
#     so the Chair/member row (Remy Giraud) spans two rows, and each of the
#     two rows holds one Vice-Chair (Douglas Body*, then Ilse Gayl*).
#  4. Region / Gender roll-up tables (shifted down one row by the insert)
#     get updated counts: Region IV 4->5, Region V 0->1, Gender M 8->9,
#     Gender F 5->6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Vice-Chair column header gains Co-Chair wording + footnote mark ---
$ws.Range("G1").Value = "Co-Chair/Vice-Chair*"

# --- 2. Existing Vice-Chair for TT-AC gets a footnote marker ---
$ws.Range("G12").Value = "Eugene Burger*"

# --- 3. Insert a row after row 13 so the TT-First-Mile entry can list two
#        Vice-Chairs (one per row) while Chair info spans both rows ---
$ws.Rows(14).Insert()

# Re-merge the Task Team category label down across the new row
$ws.Range("A11:A13").UnMerge()
$ws.Range("A11:A14").Merge()

# First new Vice-Chair goes on the existing (now row 13) line
$ws.Range("G13").Value = "Douglas Body*"
$ws.Range("H13").Value = "V"
$ws.Range("I13").Value = "Australia"
$ws.Range("J13").Value = "M"

# Second new Vice-Chair occupies the freshly inserted row 14
$ws.Range("G14").Value = "Ilse Gayl*"
$ws.Range("H14").Value = "IV"
$ws.Range("I14").Value = "USA"
$ws.Range("J14").Value = "F"

# Merge the Chair / member columns across the new two-row block so they
# continue to read as a single entry
$ws.Range("B13:B14").Merge()
$ws.Range("C13:C14").Merge()
$ws.Range("D13:D14").Merge()
$ws.Range("E13:E14").Merge()
$ws.Range("F13:F14").Merge()

# --- 4. Update the Region / Gender summary tables (now one row lower) ---
$ws.Range("B21").Value = 5   # Region IV count
$ws.Range("B22").Value = 1   # Region V count
$ws.Range("B27").Value = 9   # Gender M count
$ws.Range("B28").Value = 6   # Gender F count
